$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string must be forced to
# Text format first, otherwise Excel auto-converts the assignment to a number
# (the source data keeps these as literal text, matching the original inlineStr cells).
$ws.Range("D2").Value = '26.488.43'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '1.618.92'
$ws.Range("E3").Value = '  +1.87%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.26'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.504'
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.27'
$ws.Range("E10").Value = '  -0.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0857'
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D12").Value = '1.844.56'
$ws.Range("E12").Value = '  +1.77%  '
$ws.Range("D13").Value = '1.617.95'
$ws.Range("E13").Value = '  +1.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.04'
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.513'
$ws.Range("E15").Value = '  -1.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.78'
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = '26.487.39'
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.88'
$ws.Range("E18").Value = '  +8.22%  '
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("E20").Value = '  +2.39%  '
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.37'
$ws.Range("E22").Value = '  +2.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.13'
$ws.Range("E23").Value = '  +1.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.13'
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.38'
$ws.Range("E25").Value = '  +0.98%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("E28").Value = '  +2.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.53'
$ws.Range("E29").Value = '  +2.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0499'
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("E32").Value = '  +0.87%  '
$ws.Range("D33").Value = '1.454.67'
$ws.Range("E33").Value = '  +8.99%  '
$ws.Range("E34").Value = '  +2.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.41'
$ws.Range("E35").Value = '  -0.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.559'
$ws.Range("E37").Value = '  -5.07%  '
$ws.Range("E38").Value = '  +0.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.836'
$ws.Range("E39").Value = '  +2.43%  '
$ws.Range("E40").Value = '  +2.34%  '
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.20'
$ws.Range("E42").Value = '  +2.70%  '
$ws.Range("D43").Value = '1.756.95'
$ws.Range("E43").Value = '  +1.92%  '
$ws.Range("E44").Value = '  -0.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.15'
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.916'
$ws.Range("E46").Value = '  -9.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.11'
$ws.Range("E47").Value = '  +3.08%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("E48").Value = '  -1.21%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.50'
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0503'
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0963'
$ws.Range("E51").Value = '  -0.96%  '
